$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.401.00'
$ws.Range('E2').Value = '  -2.35%  '

$ws.Range('D3').Value = '1.838.34'
$ws.Range('E3').Value = '  -2.15%  '

$ws.Range('E4').Value = '  +0.19%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '259.84'
$ws.Range('E5').Value = '  -7.84%  '

$ws.Range('E6').Value = '  +0.16%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5192'
$ws.Range('E7').Value = '  -1.46%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3220'
$ws.Range('E8').Value = '  -8.74%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06726'
$ws.Range('E9').Value = '  -4.29%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.80'
$ws.Range('E10').Value = '  -7.62%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7674'
$ws.Range('E11').Value = '  -6.04%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07685'
$ws.Range('E12').Value = '  -1.22%  '

$ws.Range('D13').Value = '1.851.50'
$ws.Range('E13').Value = '  -1.36%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.99'
$ws.Range('E14').Value = '  -1.71%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.013'
$ws.Range('E15').Value = '  -3.77%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  +0.31%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.09'
$ws.Range('E17').Value = '  -3.40%  '

$ws.Range('E18').Value = '  +0.01%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007850'
$ws.Range('E19').Value = '  -3.99%  '

$ws.Range('D20').Value = '26.432.22'
$ws.Range('E20').Value = '  -2.32%  '

$ws.Range('D21').Value = '2.083.78'
$ws.Range('E21').Value = '  -0.94%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.527'
$ws.Range('E22').Value = '  -4.98%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.446'
$ws.Range('E23').Value = '  -7.24%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.899'
$ws.Range('E24').Value = '  -5.27%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.323'
$ws.Range('E25').Value = '  -2.47%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '144.39'
$ws.Range('E26').Value = '  -1.23%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.664'
$ws.Range('E27').Value = '  -1.03%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.89'
$ws.Range('E28').Value = '  -3.88%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '110.97'
$ws.Range('E29').Value = '  -1.83%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.167'
$ws.Range('E30').Value = '  -4.88%  '

$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08748'
$ws.Range('E31').Value = '  -1.60%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.119'
$ws.Range('E32').Value = '  -5.89%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04832'
$ws.Range('E33').Value = '  -1.26%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.129'
$ws.Range('E34').Value = '  -3.86%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.857'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6814'
$ws.Range('E36').Value = '  -8.13%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.101'
$ws.Range('E37').Value = '  -5.83%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01780'
$ws.Range('E38').Value = '  -5.36%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.209'
$ws.Range('E39').Value = '  -8.32%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4909'
$ws.Range('E40').Value = '  -7.34%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '111.42'
$ws.Range('E41').Value = '  -4.87%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8925'
$ws.Range('E42').Value = '  -8.97%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.158'
$ws.Range('E43').Value = '  -2.46%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.21%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.718'
$ws.Range('E45').Value = '  -5.83%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4192'
$ws.Range('E46').Value = '  -8.77%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.111'
$ws.Range('E47').Value = '  -3.57%  '

$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1257'
$ws.Range('E48').Value = '  -8.16%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05881'
$ws.Range('E49').Value = '  -1.04%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '35.33'
$ws.Range('E50').Value = '  -3.82%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '59.16'
$ws.Range('E51').Value = '  -4.27%  '
